$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.382.85'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').Value = '3.505.58'
$ws.Range('E3').Value = '  -3.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '200.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '551.24'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.62%  '
$ws.Range('D7').Value = '3.499.24'
$ws.Range('E7').Value = '  -3.91%  '
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '63.61'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +10.83%  '
$ws.Range('E12').Value = '  -7.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -8.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.83'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '4.059.94'
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').Value = '3.494.75'
$ws.Range('E16').Value = '  -3.80%  '
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.34'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.58%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '67.156.40'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.77'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.52%  '
$ws.Range('E21').Value = '  -5.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '391.64'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.16'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.80%  '
$ws.Range('E24').Value = '  -5.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.46'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.86'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.22'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.60%  '
$ws.Range('E28').Value = '  -5.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.81'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.95'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '676.88'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.98'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -14.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.74'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.87'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('E35').Value = '  -7.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.54'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -10.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.398'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.33%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.131'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.79%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.071.66'
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E42').Value = '  -4.95%  '
$ws.Range('D43').Value = '0.0₃0672'
$ws.Range('E43').Value = '  -15.85%  '
$ws.Range('E44').Value = '  -12.35%  '
$ws.Range('E45').Value = '  +6.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0398'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.127'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '137.27'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.24'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.90'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -7.61%  '
